# Final Examination Green Light Form -- address/town/postal/time updates.
$d = $word.ActiveDocument

# 1) Address: "Van Limburg Stirumstraat 228" -> "Bosbouw 9" (spell-checked run)
$r = $d.Content
$r.Find.Execute("Van Limburg Stirumstraat 228")
[void]$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="14571EAC" w14:textId="7269CA0F" w:rsidR="0099748E" w:rsidRPr="002E319C" w:rsidRDefault="0099748E" w:rsidP="00F074D2"><w:pPr><w:spacing w:line="20" w:lineRule="atLeast"/><w:outlineLvl w:val="0"/><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve"> </w:t></w:r><w:permStart w:id="19150507" w:edGrp="everyone"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Bosbouw</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> 9</w:t></w:r><w:r><w:t xml:space="preserve">  </w:t></w:r><w:permEnd w:id="19150507"/></w:p>')

# 2) Postal code: remove "2515PS" value, leaving blank (3 spaces)
$r = $d.Content
$r.Find.Execute("2515PS")
[void]$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="5D8B058D" w14:textId="6D0E5FE0" w:rsidR="0099748E" w:rsidRPr="002E319C" w:rsidRDefault="0099748E" w:rsidP="00F074D2"><w:pPr><w:spacing w:line="20" w:lineRule="atLeast"/><w:outlineLvl w:val="0"/><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve"> </w:t></w:r><w:permStart w:id="590641938" w:edGrp="everyone"/><w:r><w:t xml:space="preserve">   </w:t></w:r><w:permEnd w:id="590641938"/></w:p>')

# 3) Town: "Den Haag" -> "Houten" (spell-checked run)
$r = $d.Content
$r.Find.Execute("Den Haag")
[void]$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="06ECDB41" w14:textId="7FD4BE2E" w:rsidR="0099748E" w:rsidRPr="002E319C" w:rsidRDefault="0099748E" w:rsidP="00F074D2"><w:pPr><w:spacing w:line="20" w:lineRule="atLeast"/><w:outlineLvl w:val="0"/><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:r><w:t xml:space="preserve"> </w:t></w:r><w:permStart w:id="574958385" w:edGrp="everyone"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Houten</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">  </w:t></w:r><w:permEnd w:id="574958385"/></w:p>')

# 4) Time: "11:00" -> "13:00"
$r = $d.Content
$r.Find.Execute("11:00")
[void]$r.InsertXML('<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="4C2695E6" w14:textId="53E12667" w:rsidR="00037874" w:rsidRPr="00522C70" w:rsidRDefault="00037874" w:rsidP="00F13CBE"><w:pPr><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="nl-NL"/></w:rPr></w:pPr><w:permStart w:id="1067017590" w:edGrp="everyone"/><w:r w:rsidRPr="00522C70"><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="nl-NL"/></w:rPr><w:t>1</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="nl-NL"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="nl-NL"/></w:rPr><w:t>:00</w:t></w:r><w:r w:rsidRPr="00522C70"><w:rPr><w:rFonts w:ascii="Tahoma" w:hAnsi="Tahoma" w:cs="Tahoma"/><w:sz w:val="18"/><w:szCs w:val="18"/><w:lang w:val="nl-NL"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:permEnd w:id="1067017590"/></w:p>')
